$wb = $excel.ActiveWorkbook

$changes = @{
    3  = 2965
    7  = 1617
    11 = 1325
    13 = 441
    14 = 335
    20 = 3042
    21 = 371
    22 = 84
    24 = 83
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $changes.Keys) {
        $ws.Cells.Item($row, 6).Value = $changes[$row]
    }
}
